$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("I2").Value = 3.7
$ws.Range("K2").Value = 1.91
$ws.Range("M2").Value = 1.11
$ws.Range("N2").Value = 6.5
$ws.Range("O2").Value = 1.53
$ws.Range("P2").Value = 2.38
$ws.Range("Q2").Value = 2.7
$ws.Range("R2").Value = 1.44
$ws.Range("U2").Value = 2.2
$ws.Range("V2").Value = 1.62
$ws.Range("AC2").Value = 6
$ws.Range("AF2").Value = 81
$ws.Range("AI2").Value = 15
$ws.Range("AL2").Value = 51
$ws.Range("AS2").Value = 301
$ws.Range("AU2").Value = 9.5
$ws.Range("AV2").Value = 81
$ws.Range("G3").Value = 2.15
$ws.Range("H3").Value = 2.88
$ws.Range("I3").Value = 3.9
$ws.Range("J3").Value = 3.1
$ws.Range("K3").Value = 1.83
$ws.Range("M3").Value = 1.14
$ws.Range("N3").Value = 5.5
$ws.Range("S3").Value = 1.67
$ws.Range("T3").Value = 2.1
$ws.Range("Y3").Value = 11
$ws.Range("Z3").Value = 21
$ws.Range("AG3").Value = 7.5
$ws.Range("AH3").Value = 17
$ws.Range("AP3").Value = 34
$ws.Range("AQ3").Value = 51
$ws.Range("AS3").Value = 351
$ws.Range("AT3").Value = 2.1
$ws.Range("AZ3").Value = 81
$ws.Range("I5").Value = 3.9
$ws.Range("M5").Value = 1.11
$ws.Range("N5").Value = 6.5
$ws.Range("X5").Value = 8
$ws.Range("Z5").Value = 17
$ws.Range("AN5").Value = 3.75
$ws.Range("BA5").Value = 151
$ws.Range("Q6").Value = 2.4
$ws.Range("R6").Value = 1.53
